$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2205.9155
$ws.Range("I15").Value = 2205.9155
$ws.Range("K15").Value = 6617.7465
$ws.Range("M15").Value = -6448.7465

$ws.Range("H86").Value = 12207.737
$ws.Range("I86").Value = 17339.7
$ws.Range("J86").Value = 6505.5557
$ws.Range("K86").Value = 17339.7
$ws.Range("L86").Value = 6505.5557
$ws.Range("M86").Value = -16216.7
$ws.Range("N86").Value = -8751.555700000001

$ws.Range("H88").Value = 419003.9
$ws.Range("J88").Value = 511123.7
$ws.Range("L88").Value = 511123.7
$ws.Range("N88").Value = -511935.7

$ws.Range("H89").Value = 12207.737
$ws.Range("I89").Value = 17339.7
$ws.Range("J89").Value = 6505.5557
$ws.Range("K89").Value = 86698.5
$ws.Range("L89").Value = 32527.7785
$ws.Range("M89").Value = -81082.5
$ws.Range("N89").Value = -43759.7785

$ws.Range("H91").Value = 419003.9
$ws.Range("J91").Value = 511123.7
$ws.Range("L91").Value = 511123.7
$ws.Range("N91").Value = -513931.7

$ws.Range("H96").Value = 1606.1818
$ws.Range("I96").Value = 1032
$ws.Range("J96").Value = 2611
$ws.Range("K96").Value = 3096
$ws.Range("L96").Value = 7833
$ws.Range("M96").Value = -1723
$ws.Range("N96").Value = -10579

$ws.Range("H106").Value = 7438.8335
$ws.Range("I106").Value = 2035.6
$ws.Range("K106").Value = 2035.6
$ws.Range("M106").Value = -1404.6

$ws.Range("H132").Value = 2453.9473
$ws.Range("I132").Value = 2423.6667
$ws.Range("K132").Value = 7271.000100000001
$ws.Range("M132").Value = -4741.000100000001

$ws.Range("H138").Value = 5665.815
$ws.Range("J138").Value = 5386.109
$ws.Range("L138").Value = 16158.327
$ws.Range("N138").Value = -26438.327

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8207852.5
$ws.Range("I32").Value = 9267076
$ws.Range("J32").Value = 36701.855
$ws.Range("K32").Value = 9267076
$ws.Range("L32").Value = 36701.855
$ws.Range("M32").Value = -9266789
$ws.Range("N32").Value = -37275.855

$ws.Range("H45").Value = 2422.7222
$ws.Range("I45").Value = 2330.1428
$ws.Range("K45").Value = 2330.1428
$ws.Range("M45").Value = -1953.1428

$ws.Range("H61").Value = 23862084
$ws.Range("I61").Value = 41670110
$ws.Range("J61").Value = 118047.11
$ws.Range("K61").Value = 41670110
$ws.Range("L61").Value = 118047.11
$ws.Range("M61").Value = -41669898
$ws.Range("N61").Value = -118471.11

$ws.Range("H74").Value = 22743476
$ws.Range("I74").Value = 125000600
$ws.Range("K74").Value = 125000600
$ws.Range("M74").Value = -124999726

$ws.Range("H77").Value = 22743476
$ws.Range("I77").Value = 125000600
$ws.Range("K77").Value = 625003000
$ws.Range("M77").Value = -624998632

$ws.Range("H88").Value = 2386.5454
$ws.Range("I88").Value = 2107.5
$ws.Range("J88").Value = 2721.4
$ws.Range("K88").Value = 2107.5
$ws.Range("L88").Value = 2721.4
$ws.Range("M88").Value = -1701.5
$ws.Range("N88").Value = -3533.4

$ws.Range("H91").Value = 2386.5454
$ws.Range("I91").Value = 2107.5
$ws.Range("J91").Value = 2721.4
$ws.Range("K91").Value = 2107.5
$ws.Range("L91").Value = 2721.4
$ws.Range("M91").Value = -703.5
$ws.Range("N91").Value = -5529.4

$ws.Range("H122").Value = 3224.3333
$ws.Range("I122").Value = 2337.625
$ws.Range("K122").Value = 7012.875
$ws.Range("M122").Value = -4562.875

$ws.Range("H132").Value = 5884.3784
$ws.Range("I132").Value = 2849.84
$ws.Range("K132").Value = 8549.52
$ws.Range("M132").Value = -6019.52

$ws.Range("H136").Value = 23862084
$ws.Range("I136").Value = 41670110
$ws.Range("J136").Value = 118047.11
$ws.Range("K136").Value = 125010330
$ws.Range("L136").Value = 354141.33
$ws.Range("M136").Value = -125007780
$ws.Range("N136").Value = -359241.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2267.5715
$ws.Range("I107").Value = 1829.0834
$ws.Range("K107").Value = 1829.0834
$ws.Range("M107").Value = 90.91660000000002

$ws.Range("H134").Value = 32139.629
$ws.Range("I134").Value = 3060
$ws.Range("K134").Value = 9180
$ws.Range("M134").Value = -6645

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1978.8
$ws.Range("I105").Value = 1754.25
$ws.Range("J105").Value = 2128.5
$ws.Range("K105").Value = 1754.25
$ws.Range("L105").Value = 2128.5
$ws.Range("M105").Value = -7.25
$ws.Range("N105").Value = -5622.5

$ws.Range("H111").Value = 70702
$ws.Range("J111").Value = 70702
$ws.Range("L111").Value = 70702
$ws.Range("N111").Value = -78882

$ws.Range("H122").Value = 2822.5862
$ws.Range("I122").Value = 1340.0588
$ws.Range("J122").Value = 4922.8335
$ws.Range("K122").Value = 4020.1764
$ws.Range("L122").Value = 14768.5005
$ws.Range("M122").Value = -1570.1764
$ws.Range("N122").Value = -19668.5005

$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 751
$ws.Range("I5").Value = 542
$ws.Range("K5").Value = 1626
$ws.Range("M5").Value = -1514

$ws.Range("H68").Value = 2117.2424
$ws.Range("I68").Value = 2364.8572
$ws.Range("J68").Value = 2050.577
$ws.Range("K68").Value = 7094.571599999999
$ws.Range("L68").Value = 6151.731000000001
$ws.Range("M68").Value = -6283.571599999999
$ws.Range("N68").Value = -7773.731000000001

$ws.Range("H71").Value = 2117.2424
$ws.Range("I71").Value = 2364.8572
$ws.Range("J71").Value = 2050.577
$ws.Range("K71").Value = 21283.7148
$ws.Range("L71").Value = 18455.193
$ws.Range("M71").Value = -17227.7148
$ws.Range("N71").Value = -26567.193

$ws.Range("H107").Value = 785.1923
$ws.Range("I107").Value = 751.8182
$ws.Range("J107").Value = 809.6667
$ws.Range("K107").Value = 2255.4546
$ws.Range("L107").Value = 2429.0001
$ws.Range("M107").Value = -335.4546
$ws.Range("N107").Value = -6269.0001

$ws.Range("H108").Value = 1499
$ws.Range("I108").Value = 1499
$ws.Range("K108").Value = 4497
$ws.Range("M108").Value = -1617

$ws.Range("H135").Value = 751
$ws.Range("I135").Value = 542
$ws.Range("K135").Value = 4878
$ws.Range("M135").Value = -2343

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9003
$ws.Range("J70").Value = 9003
$ws.Range("L70").Value = 9003
$ws.Range("N70").Value = -9543

$ws.Range("H73").Value = 9003
$ws.Range("J73").Value = 9003
$ws.Range("L73").Value = 9003
$ws.Range("N73").Value = -10875

$ws.Range("H97").Value = 5068.3335
$ws.Range("I97").Value = 5082
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 5082
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -4586
$ws.Range("N97").Value = -5992

$ws.Range("H122").Value = 22709.375
$ws.Range("I122").Value = 44000
$ws.Range("K122").Value = 132000
$ws.Range("M122").Value = -129550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41908.54
$ws.Range("I7").Value = 2207
$ws.Range("J7").Value = 96047
$ws.Range("K7").Value = 2207
$ws.Range("L7").Value = 96047
$ws.Range("M7").Value = -2095
$ws.Range("N7").Value = -96271

$ws.Range("H40").Value = 3257.923
$ws.Range("I40").Value = 3112.75
$ws.Range("K40").Value = 3112.75
$ws.Range("M40").Value = -2976.75

$ws.Range("H55").Value = 58823956
$ws.Range("I55").Value = 76923496
$ws.Range("J55").Value = 434.25
$ws.Range("K55").Value = 76923496
$ws.Range("L55").Value = 434.25
$ws.Range("M55").Value = -76923323
$ws.Range("N55").Value = -780.25

$ws.Range("H68").Value = 1091.6
$ws.Range("J68").Value = 897
$ws.Range("L68").Value = 897
$ws.Range("N68").Value = -2395

$ws.Range("H71").Value = 1091.6
$ws.Range("J71").Value = 897
$ws.Range("L71").Value = 4485
$ws.Range("N71").Value = -11973

$ws.Range("H74").Value = 15818998
$ws.Range("I74").Value = 31525000
$ws.Range("J74").Value = 112995
$ws.Range("K74").Value = 31525000
$ws.Range("L74").Value = 112995
$ws.Range("M74").Value = -31524002
$ws.Range("N74").Value = -114991

$ws.Range("H77").Value = 15818998
$ws.Range("I77").Value = 31525000
$ws.Range("J77").Value = 112995
$ws.Range("K77").Value = 94575000
$ws.Range("L77").Value = 338985
$ws.Range("M77").Value = -94570008
$ws.Range("N77").Value = -348969

$ws.Range("H93").Value = 32259778
$ws.Range("I93").Value = 38463330
$ws.Range("K93").Value = 38463330
$ws.Range("M93").Value = -38462082

$ws.Range("H126").Value = 41908.54
$ws.Range("I126").Value = 2207
$ws.Range("J126").Value = 96047
$ws.Range("K126").Value = 6621
$ws.Range("L126").Value = 288141
$ws.Range("M126").Value = -4151
$ws.Range("N126").Value = -293081

$ws.Range("H132").Value = 45206.04
$ws.Range("I132").Value = 46974.32
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 140922.96
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -138392.96
$ws.Range("N132").Value = -8057

$ws.Range("H136").Value = 44155.688
$ws.Range("I136").Value = 6880.0625
$ws.Range("J136").Value = 81431.31
$ws.Range("K136").Value = 20640.1875
$ws.Range("L136").Value = 244293.93
$ws.Range("M136").Value = -18090.1875
$ws.Range("N136").Value = -249393.93

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 58750
$ws.Range("I92").Value = 55000
$ws.Range("J92").Value = 60000
$ws.Range("K92").Value = 55000
$ws.Range("L92").Value = 60000
$ws.Range("M92").Value = -52504
$ws.Range("N92").Value = -64992

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
